# Updated cryptos list on Sun Jun 23 07:49:01 UTC 2024 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) values for the
# crypto rows in the sheet. The Price column holds plain text that often
# *looks* numeric (e.g. "592.21", "64.385.23", "0.0000117"); assigning such a
# string straight to .Value would make Excel auto-coerce it into a genuine
# floating point number (losing the original text formatting/precision and
# adding a number-format style). To keep these as plain text - exactly like
# the original workbook stored them - we briefly force a text number format
# before writing the value and then clear the formatting again so the cell's
# style stays at its original (default) state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($address, $text) {
    $cell = $ws.Range($address)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

function Set-PlainValue($address, $text) {
    $ws.Range($address).Value = $text
}

# row -> Price (D), Volume(1h) (E)
$updates = @(
    @{ Row = 2;  D = '64.385.23';  E = '  -0.12%  ' },
    @{ Row = 3;  D = '3.509.62' },
    @{ Row = 4;  E = '  +0.01%  ' },
    @{ Row = 5;  D = '592.21' },
    @{ Row = 6;  D = '134.74';    E = '  -0.26%  ' },
    @{ Row = 7;  E = '  +0.00%  ' },
    @{ Row = 8;  E = '  +0.14%  ' },
    @{ Row = 9;  D = '7.61';      E = '  +5.67%  ' },
    @{ Row = 10; E = '  +0.62%  ' },
    @{ Row = 11; D = '0.389';     E = '  +3.46%  ' },
    @{ Row = 12; D = '4.107.79';  E = '  +0.39%  ' },
    @{ Row = 13; E = '  +1.18%  ' },
    @{ Row = 14; E = '  +0.51%  ' },
    @{ Row = 15; D = '3.510.47';  E = '  +0.43%  ' },
    @{ Row = 16; D = '25.77';     E = '  +1.79%  ' },
    @{ Row = 17; D = '64.379.43'; E = '  -0.10%  ' },
    @{ Row = 18; E = '  -0.11%  ' },
    @{ Row = 19; D = '13.65';     E = '  -0.88%  ' },
    @{ Row = 20; D = '5.76';      E = '  +2.05%  ' },
    @{ Row = 21; D = '394.41';    E = '  +2.18%  ' },
    @{ Row = 22; E = '  +2.03%  ' },
    @{ Row = 23; D = '3.649.56';  E = '  +0.41%  ' },
    @{ Row = 24; D = '74.63';     E = '  +0.76%  ' },
    @{ Row = 25; E = '  -0.01%  ' },
    @{ Row = 26; D = '5.76';      E = '  +0.63%  ' },
    @{ Row = 27; D = '0.0000117'; E = '  +3.24%  ' },
    @{ Row = 28; E = '  -0.02%  ' },
    @{ Row = 29; D = '7.42';      E = '  -0.49%  ' },
    @{ Row = 31; D = '8.26';      E = '  +0.39%  ' },
    @{ Row = 32; E = '  -3.93%  ' },
    @{ Row = 33; E = '  +7.05%  ' },
    @{ Row = 34; D = '3.538.34';  E = '  +0.53%  ' },
    @{ Row = 35; E = '  +0.02%  ' },
    @{ Row = 36; D = '23.37';     E = '  -0.45%  ' },
    @{ Row = 37; D = '5.38';      E = '  +1.52%  ' },
    @{ Row = 38; D = '6.96';      E = '  +1.72%  ' },
    @{ Row = 39; E = '  +1.03%  ' },
    @{ Row = 40; D = '167.05';    E = '  +2.83%  ' },
    @{ Row = 41; D = '0.0788';    E = '  +0.96%  ' },
    @{ Row = 42; D = '0.812';     E = '  +0.95%  ' },
    @{ Row = 43; E = '  +0.03%  ' },
    @{ Row = 44; E = '  +0.99%  ' },
    @{ Row = 45; D = '24.84';     E = '  -4.12%  ' },
    @{ Row = 46; E = '  +0.42%  ' },
    @{ Row = 47; E = '  -3.08%  ' },
    @{ Row = 48; D = '6.81';      E = '  +0.58%  ' },
    @{ Row = 49; E = '  +0.58%  ' },
    @{ Row = 50; D = '2.376.76';  E = '  -4.11%  ' },
    @{ Row = 51; D = '0.0261';    E = '  +0.33%  ' }
)

foreach ($u in $updates) {
    if ($u.ContainsKey('D')) {
        Set-TextValue ("D" + $u.Row) $u.D
    }
    if ($u.ContainsKey('E')) {
        Set-PlainValue ("E" + $u.Row) $u.E
    }
}
